# JCSantos - Casi OK. Antes de informe
#
# Adds two new Name/Value settings rows to the "Settings" sheet:
#   Row 19: CorreoHtmlSinElementos / Extra\cuerpo3.html
#   Row 20: AsuntoSinTransacciones / Sin correos que procesar
# and nudges the sheet's used range / selection to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New "no transactions / empty email" settings, appended right after the
# existing AsuntoCorreoOK row (row 18).
$ws.Range("A19").Value = "CorreoHtmlSinElementos"
$ws.Range("B19").Value = "Extra\cuerpo3.html"

$ws.Range("A20").Value = "AsuntoSinTransacciones"
$ws.Range("B20").Value = "Sin correos que procesar"

# The sheet's formatted area grows by one more (still empty) row.
$ws.Rows.Item(1001).RowHeight = 14.25

# Leave the selection where it ended up when the workbook was last saved.
$ws.Range("A16").Select()
